$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    12 = 5.049999999999998
    27 = 5.981400000000007
    32 = 6.701299999999998
    36 = 9.119200000000001
    38 = 5.040100000000001
    46 = 6.478899999999999
    54 = 4.650599999999995
    55 = 5.437699999999998
    56 = 4.851199999999997
    67 = 5.449599999999996
    69 = 5.157399999999996
    72 = 4.871300000000003
    83 = 5.6107
    86 = 4.7746
    91 = 5.4212
    93 = 5.794200000000003
    99 = 4.742399999999997
}

foreach ($row in $updates.Keys) {
    $ws.Range("B$row").Value = $updates[$row]
}
